# The document contains the tag/id/tag sequence "<id>p161v_1</id>" split
# across three separate runs (the surrounding "<id>" / "</id>" tag runs use
# Courier New / color 7f6000 / sz 18, while the middle "p161v_1" run uses
# plain default formatting). The edit collapses all three runs into a
# single run (formatted like the original "<id>" run) whose text is the
# full, concatenated "<id>p161v_1</id>" string.

$d = $word.ActiveDocument

# Locate the run of text spanning the three original runs.
$rng = $d.Content
$found = $rng.Find.Execute("<id>p161v_1</id>", $false, $false, $false, `
                            $false, $false, $true, 1, $false, "", 0)

if (-not $found) {
    throw "Could not find '<id>p161v_1</id>' in the document"
}

# Re-assigning Range.Text with the *same* string is treated as a no-op by
# the host, so nothing would be merged. Swap in a throwaway placeholder
# first -- that collapses the whole (currently 3-run) range down to a
# single run carrying the first run's formatting -- then put the real
# text back into that same (now single-run) range.
$start = $rng.Start
$end = $rng.End
$rng.Text = "TEMP_PLACEHOLDER"

$rng2 = $d.Range($start, $rng.End)
$rng2.Text = "<id>p161v_1</id>"
